$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# --- "Status" rolls from "Ready for handoff" to "Handed back: in sync with en-US"
# (same text everywhere it's used: the Overview rollup columns as well as each
# language sheet's own Status column).
$newStatus = "Handed back: in sync with en-US"
$ovw.Range("B2").Value = $newStatus
$ovw.Range("C2").Value = $newStatus
$ovw.Range("B3").Value = $newStatus
$ovw.Range("C3").Value = $newStatus
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# --- Latest Handback DateTime (column H) now carries a real handback timestamp
# per language instead of the zero-date placeholder.
$zh.Range("H2").Value = "2016-03-21 18:32:29"
$zh.Range("H3").Value = "2016-03-21 18:32:29"
$de.Range("H2").Value = "2016-03-21 18:32:36"
$de.Range("H3").Value = "2016-03-21 18:32:36"

# --- New "Latest Target File" (F) / "Latest Handback File" (G) hyperlinked
# columns. F mirrors the source-file hyperlink (column A); G mirrors the
# handoff-file hyperlink (column D) -- same display text & target URL, for
# every data row on both language sheets.

function Add-ReportHyperlink($sheet, $cellRef, $address, $displayText) {
    $sheet.Range($cellRef).Style = "HyperLink"
    $sheet.Hyperlinks.Add($sheet.Range($cellRef), $address, "", "", $displayText) | Out-Null
}

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7863e071f6e9e7329af1684e980b9587c4b23bd9/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a702d86c53ad2e94556664c04fe1f5bf8c1879b9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/930540005c2cae433e47edfaef9eba1844babd86/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# zh-cn, rows 2 & 3
Add-ReportHyperlink $zh "F2" $mdUrl "a.md"
Add-ReportHyperlink $zh "G2" $zhXlfUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Add-ReportHyperlink $zh "F3" $mdUrl "a.md"
Add-ReportHyperlink $zh "G3" $zhXlfUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# de-de, rows 2 & 3
Add-ReportHyperlink $de "F2" $mdUrl "a.md"
Add-ReportHyperlink $de "G2" $deXlfUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Add-ReportHyperlink $de "F3" $mdUrl "a.md"
Add-ReportHyperlink $de "G3" $deXlfUrl "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
